$d = $word.ActiveDocument

# 1) Update the letter date.
$d.Content.Find.Execute("June 30, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "October 30, 2025", 2) | Out-Null

# 2) Locate the paragraph that confirms the Folio number(s).
$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Hereby, the legal heirs*") {
        $para = $p
    }
}

$pStart = $para.Range.Start
$pEnd = $para.Range.End

# Collect the SUBJECT field and the (now-unused) second "Folio2" DOCPROPERTY
# field that live inside this paragraph - these are the merge fields that
# used to build up a "Folio2 & Folio1" list. Leave the Folio1 field alone.
$toRemove = @()
foreach ($f in $d.Fields) {
    if ($f.Code.Start -ge $pStart -and $f.Code.Start -lt $pEnd) {
        if ($f.Code.Text -match "SUBJECT" -or $f.Code.Text -match "Folio2") {
            $toRemove += $f
        }
    }
}

# Delete from the highest start position down so earlier offsets/fields
# remain valid while we work.
$toRemove = $toRemove | Sort-Object -Property {$_.Code.Start} -Descending
foreach ($f in $toRemove) {
    $f.Delete()
}

# 3) Drop the now-orphaned " & " joiner that used to sit between the two
# folio fields.
$para.Range.Find.Execute(" & ", $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 2) | Out-Null

# 4) The paragraph now references a single folio, so "Folio Nos" becomes
# "Folio No".
$para.Range.Find.Execute("Folio Nos ", $true, $false, $false, $false, $false,
                          $true, 1, $false, "Folio No ", 2) | Out-Null
